$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new 7th table column "LeaveEmpty" to Table3 (was A1:F5 -> A1:G5)
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null
$ws.Range("G1").Value2 = "LeaveEmpty"

# Fill in the previously-empty Allergens cell for the 4th menu item
$ws.Range("C4").Value2 = "No known priority allergens"
$ws.Range("C4").WrapText = $true

# Update selection to match the saved workbook state
$ws.Range("E12").Select()
